$wb = $excel.ActiveWorkbook

# --- Reorder worksheets: "review_info" becomes the first tab, "hotel_info" the second ---
$wsHotel = $wb.Worksheets.Item("hotel_info")
$wsReview = $wb.Worksheets.Item("review_info")
$wsReview.Move($wsHotel)

# --- hotel_info: insert a new "State" column between "Hotel_Name" and "City" ---
$wsHotel = $wb.Worksheets.Item("hotel_info")
$wsHotel.Columns.Item(3).Insert()
$wsHotel.Range("C1").Value = "State"
$wsHotel.Range("C2").Value = "Louisiana"
